$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the shared-string header labels -------------------------------
# "_old" -> "_FV2310" and "_new" -> "_FV2404" for the 10 paired header cells
# in row 1 (columns A-J are the *_old set, L-U are the *_new set; column K
# is the unchanged "diff" header).
$oldSuffixHeaders = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldSuffixHeaders.Count; $i++) {
    $base = $oldSuffixHeaders[$i]
    $ws.Cells.Item(1, $i + 1).Value = "$($base)_FV2310"
    $ws.Cells.Item(1, $i + 12).Value = "$($base)_FV2404"
}

# --- Freeze the header row -------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Turn the used range into a real Excel Table ---------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U80"), $false, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
